$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestData = $wb.Worksheets.Item("TestData")

# Toggle Runmode values
$wsTestCases.Range("B3").Value = "Y"
$wsTestData.Range("A8").Value = "N"

# Update active sheet / selection state: TestCases becomes the active tab,
# selection on TestData moves to A8.
$wsTestData.Range("A8").Select()
$wsTestCases.Activate()
$wsTestCases.Range("B4").Select()
